# Apply updated statistics after rerunning analysis with corrected recording dates.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: normality
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("normality")

$ws.Range("C3").Value = 0.9611
$ws.Range("D3").Value = 0.0325

$ws.Range("C5").Value = 0.9734
$ws.Range("D5").Value = 0.1532

$ws.Range("C7").Value = 0.9599
$ws.Range("D7").Value = 0.0281

$ws.Range("C9").Value = 0.9572000000000001
$ws.Range("D9").Value = 0.02

$ws.Range("D11").Value = 0.0033

# ---------------------------------------------------------------------------
# Sheet: equal_var
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("equal_var")

$ws.Range("C3").Value = 0.009900000000000001
$ws.Range("D3").Value = 0.9208

$ws.Range("C4").Value = 1.8035
$ws.Range("D4").Value = 0.1811

$ws.Range("C5").Value = 0.0216
$ws.Range("D5").Value = 0.8832

$ws.Range("C6").Value = 0.1177
$ws.Range("D6").Value = 0.7319

$ws.Range("C7").Value = 2.089
$ws.Range("D7").Value = 0.1502

# ---------------------------------------------------------------------------
# Sheet: mixed_anova
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mixed_anova")

$ws.Range("C3").Value = 5.3028
$ws.Range("E3").Value = 167
$ws.Range("F3").Value = 5.3028
$ws.Range("G3").Value = 0.0551
$ws.Range("H3").Value = 0.8147
$ws.Range("J3").Value = 0.0003

$ws.Range("C4").Value = 4063.2874
$ws.Range("E4").Value = 668
$ws.Range("F4").Value = 1015.8218
$ws.Range("G4").Value = 68.0325
$ws.Range("J4").Value = 0.2895
$ws.Range("K4").Value = 0.4913
$ws.Range("M4").Value = 0.1004

$ws.Range("C5").Value = 22.3087
$ws.Range("E5").Value = 668
$ws.Range("F5").Value = 5.5772
$ws.Range("G5").Value = 0.3735
$ws.Range("H5").Value = 0.8276
$ws.Range("J5").Value = 0.0022

# ---------------------------------------------------------------------------
# Sheet: pairwise_ttests
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("pairwise_ttests")

$ws.Range("H3").Value = -12.0498
$ws.Range("I3").Value = 168
$ws.Range("L3").Value = "2.32e+21"
$ws.Range("M3").Value = -0.8262

$ws.Range("H4").Value = -13.5737
$ws.Range("I4").Value = 168
$ws.Range("L4").Value = "4.067e+25"
$ws.Range("M4").Value = -1.2263

$ws.Range("H5").Value = -10.6019
$ws.Range("I5").Value = 168
$ws.Range("L5").Value = "2.359e+17"
$ws.Range("M5").Value = -1.097

$ws.Range("H6").Value = -7.5767
$ws.Range("I6").Value = 168
$ws.Range("L6").Value = "3.177e+09"
$ws.Range("M6").Value = -0.8264

$ws.Range("H7").Value = -8.192399999999999
$ws.Range("I7").Value = 168
$ws.Range("L7").Value = "1.05e+11"
$ws.Range("M7").Value = -0.4373

$ws.Range("H8").Value = -4.5678
$ws.Range("I8").Value = 168
$ws.Range("L8").Value = "1372.761"
$ws.Range("M8").Value = -0.3596

$ws.Range("H9").Value = -0.8544
$ws.Range("I9").Value = 168
$ws.Range("K9").Value = 0.3941
$ws.Range("L9").Value = "0.123"
$ws.Range("M9").Value = -0.0784

$ws.Range("H10").Value = 1.1264
$ws.Range("I10").Value = 168
$ws.Range("K10").Value = 0.2616
$ws.Range("L10").Value = "0.16"
$ws.Range("M10").Value = 0.0552

$ws.Range("H11").Value = 5.212
$ws.Range("I11").Value = 168
$ws.Range("L11").Value = "2.076e+04"
$ws.Range("M11").Value = 0.3334

$ws.Range("H12").Value = 7.6716
$ws.Range("I12").Value = 168
$ws.Range("L12").Value = "5.4e+09"
$ws.Range("M12").Value = 0.2659

$ws.Range("H13").Value = -0.2426
$ws.Range("I13").Value = 158.5964
$ws.Range("K13").Value = 0.8086
$ws.Range("L13").Value = "0.174"
$ws.Range("M13").Value = -0.0367

$ws.Range("H14").Value = -0.5328000000000001
$ws.Range("I14").Value = 147.3162
$ws.Range("K14").Value = 0.595
$ws.Range("L14").Value = "0.193"
$ws.Range("M14").Value = -0.08260000000000001

$ws.Range("H15").Value = 0.3916
$ws.Range("I15").Value = 160.6074
$ws.Range("K15").Value = 0.6958
$ws.Range("L15").Value = "0.182"
$ws.Range("M15").Value = 0.0588

$ws.Range("H16").Value = -0.0623
$ws.Range("I16").Value = 142.385
$ws.Range("K16").Value = 0.9504
$ws.Range("L16").Value = "0.17"
$ws.Range("M16").Value = -0.0098

$ws.Range("H17").Value = -0.6867
$ws.Range("I17").Value = 146.3628
$ws.Range("K17").Value = 0.4934
$ws.Range("L17").Value = "0.211"
$ws.Range("M17").Value = -0.1067

$ws.Range("H18").Value = -0.052
$ws.Range("I18").Value = 159.6463
$ws.Range("K18").Value = 0.9586
$ws.Range("L18").Value = "0.17"
$ws.Range("M18").Value = -0.0078
